$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,4,5,6,7,8) got cyclically rotated: the contents of one
# row replace the contents of another (dates, variety, quality, volume,
# min/max/avg price, origin and $/kg all move together), while columns
# A,B,C,E,F,G,H,I,J,Q,T stay constant across all rows (they are not part
# of the diff). Row 3 and row 9 are untouched.
#
# Mapping (new row content <= old row content):
#   2 <= 4
#   4 <= 8
#   5 <= 2
#   6 <= 5
#   7 <= 6
#   8 <= 7

# Snapshot original values for the columns that change, before overwriting.
$cols = @("D", "K", "L", "M", "N", "O", "P", "R", "S")
$rows = @(2, 4, 5, 6, 7, 8)

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

$mapping = @{
    2 = 4
    4 = 8
    5 = 2
    6 = 5
    7 = 6
    8 = 7
}

foreach ($newRow in $rows) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value2 = $src[$c]
    }
}
